$d = $word.ActiveDocument

# The end of the document used to contain, right after the bibliography's
# closing sentence ("Artigos de revistas..."):
#   - an empty paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# That whole block (the empty paragraph plus the two site-footer
# paragraphs) is removed, while the surrounding paragraphs are left
# untouched.

$jupiterIndex = -1
$copyrightIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*Ver no Jupiter*") {
        $jupiterIndex = $i
    }
    elseif ($text -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -gt 0 -and $copyrightIndex -eq ($jupiterIndex + 1)) {
    # The blank paragraph immediately preceding "Ver no Jupiter..." is
    # removed together with the two text paragraphs.
    $blankIndex = $jupiterIndex - 1

    $startPara = $d.Paragraphs.Item($blankIndex)
    $endPara = $d.Paragraphs.Item($copyrightIndex)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
